$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new presentation topics into column C for rows 15-19.
# Values are written in the order that yields the same shared-string
# table ordering as the canonical workbook (new strings appended in
# first-seen order: Farmbot, AutoVehicle, SmartHive, LabFlow..., Solar).
$ws.Range("C15").Value = "Farmbot"
$ws.Range("C17").Value = "AutoVehicle"
$ws.Range("C18").Value = "SmartHive"
$ws.Range("C16").Value = "LabFlow - Presenting on Tuesday instead of Saturday"
$ws.Range("C19").Value = "Solar"

# Update selected cell to match the new active cell in the sheet view
$ws.Range("C19").Select()
